# Update "想去人数" (column F) counts for the events that changed, on both
# the "展览" sheet and the "全部类型" sheet (which duplicates the same rows).

$wb = $excel.ActiveWorkbook

# row => [old, new] values for column F, as they appear on the "展览" sheet
$changesByName = @{
    "苏州·X-party 国漫游戏嘉年华03"                                   = 572
    "苏州·绘时国乙1.0-秩序之外"                                        = 494
    "苏州·首届Redamancy动漫游戏嘉年华"                                  = 1245
    "苏州·梦幻岛 国乙主题文化展（日夜场） 梦幻岛之约3.0"                      = 1072
    "昆山·第十二届理想乡动漫游戏展"                                      = 14038
    "苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区"        = 15421
    "苏州·萌动次元动漫游戏嘉年毕"                                       = 36
    "昆山·第十二届理想乡动漫游戏展嘉宾漠小然签售会"                            = 23
    "苏州·苏州湾动漫游戏嘉年华"                                        = 73
    "常熟·CDW·动漫展03"                                             = 1192
    "苏州·OCG国潮动漫游戏嘉年华"                                       = 6001
    "苏州·OCG国潮动漫游戏嘉年华阿杰内场"                                  = 956
    "苏州·YoungComic动漫嘉年华"                                      = 1079
    "【会员购严选】苏州·Come in joy动漫国潮文化节"                         = 5522
    "苏州·归离之缘原神only展"                                         = 69
    "苏州·白日梦想7.20全职猎人ONLY展"                                  = 99
    "苏州·萤火国潮文化节动漫品牌博览会"                                    = 453
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # find the last used row on column C (name column)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

    for ($r = 2; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $changesByName.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value2 = $changesByName[$name]
        }
    }
}
